$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness (column C) values for rows 2-146
# These represent a re-run of the genetic algorithm logging (different seed/test),
# producing a different best-fitness-so-far convergence curve.
$fitnessValues = New-Object 'object[,]' 145,1
$fitnessValues[0,0] = 10941
$fitnessValues[1,0] = 10941
$fitnessValues[2,0] = 10941
$fitnessValues[3,0] = 10073
$fitnessValues[4,0] = 10073
$fitnessValues[5,0] = 10073
$fitnessValues[6,0] = 10073
$fitnessValues[7,0] = 9991
$fitnessValues[8,0] = 9991
$fitnessValues[9,0] = 9713
$fitnessValues[10,0] = 9713
$fitnessValues[11,0] = 9713
$fitnessValues[12,0] = 9713
$fitnessValues[13,0] = 9713
$fitnessValues[14,0] = 9713
$fitnessValues[15,0] = 9711
$fitnessValues[16,0] = 9711
$fitnessValues[17,0] = 9711
$fitnessValues[18,0] = 9711
$fitnessValues[19,0] = 9711
$fitnessValues[20,0] = 9711
$fitnessValues[21,0] = 9647
$fitnessValues[22,0] = 8626
$fitnessValues[23,0] = 8626
$fitnessValues[24,0] = 8626
$fitnessValues[25,0] = 7667
$fitnessValues[26,0] = 7667
$fitnessValues[27,0] = 7667
$fitnessValues[28,0] = 7667
$fitnessValues[29,0] = 7667
$fitnessValues[30,0] = 7667
$fitnessValues[31,0] = 7667
$fitnessValues[32,0] = 7660
$fitnessValues[33,0] = 7660
$fitnessValues[34,0] = 7660
$fitnessValues[35,0] = 7660
$fitnessValues[36,0] = 7660
$fitnessValues[37,0] = 7660
$fitnessValues[38,0] = 7660
$fitnessValues[39,0] = 7660
$fitnessValues[40,0] = 7534
$fitnessValues[41,0] = 7534
$fitnessValues[42,0] = 7534
$fitnessValues[43,0] = 7534
$fitnessValues[44,0] = 7534
$fitnessValues[45,0] = 7534
$fitnessValues[46,0] = 7534
$fitnessValues[47,0] = 7534
$fitnessValues[48,0] = 7534
$fitnessValues[49,0] = 7534
$fitnessValues[50,0] = 7534
$fitnessValues[51,0] = 7534
$fitnessValues[52,0] = 7534
$fitnessValues[53,0] = 7534
$fitnessValues[54,0] = 7534
$fitnessValues[55,0] = 7534
$fitnessValues[56,0] = 7534
$fitnessValues[57,0] = 7534
$fitnessValues[58,0] = 7534
$fitnessValues[59,0] = 7534
$fitnessValues[60,0] = 7534
$fitnessValues[61,0] = 7534
$fitnessValues[62,0] = 7534
$fitnessValues[63,0] = 7534
$fitnessValues[64,0] = 7534
$fitnessValues[65,0] = 7534
$fitnessValues[66,0] = 7345
$fitnessValues[67,0] = 7345
$fitnessValues[68,0] = 7345
$fitnessValues[69,0] = 7345
$fitnessValues[70,0] = 7345
$fitnessValues[71,0] = 7343
$fitnessValues[72,0] = 7343
$fitnessValues[73,0] = 7343
$fitnessValues[74,0] = 7343
$fitnessValues[75,0] = 7343
$fitnessValues[76,0] = 7343
$fitnessValues[77,0] = 7343
$fitnessValues[78,0] = 7343
$fitnessValues[79,0] = 7343
$fitnessValues[80,0] = 7343
$fitnessValues[81,0] = 7312
$fitnessValues[82,0] = 7312
$fitnessValues[83,0] = 7312
$fitnessValues[84,0] = 7312
$fitnessValues[85,0] = 7312
$fitnessValues[86,0] = 7312
$fitnessValues[87,0] = 7312
$fitnessValues[88,0] = 7312
$fitnessValues[89,0] = 7312
$fitnessValues[90,0] = 7312
$fitnessValues[91,0] = 7312
$fitnessValues[92,0] = 7312
$fitnessValues[93,0] = 7312
$fitnessValues[94,0] = 7312
$fitnessValues[95,0] = 7312
$fitnessValues[96,0] = 7312
$fitnessValues[97,0] = 7312
$fitnessValues[98,0] = 7312
$fitnessValues[99,0] = 7312
$fitnessValues[100,0] = 7312
$fitnessValues[101,0] = 7312
$fitnessValues[102,0] = 7312
$fitnessValues[103,0] = 7312
$fitnessValues[104,0] = 7312
$fitnessValues[105,0] = 7312
$fitnessValues[106,0] = 7312
$fitnessValues[107,0] = 7312
$fitnessValues[108,0] = 7312
$fitnessValues[109,0] = 7312
$fitnessValues[110,0] = 7312
$fitnessValues[111,0] = 7312
$fitnessValues[112,0] = 7312
$fitnessValues[113,0] = 7312
$fitnessValues[114,0] = 7312
$fitnessValues[115,0] = 7312
$fitnessValues[116,0] = 7312
$fitnessValues[117,0] = 7312
$fitnessValues[118,0] = 7312
$fitnessValues[119,0] = 7312
$fitnessValues[120,0] = 7312
$fitnessValues[121,0] = 7312
$fitnessValues[122,0] = 7312
$fitnessValues[123,0] = 7310
$fitnessValues[124,0] = 7310
$fitnessValues[125,0] = 7310
$fitnessValues[126,0] = 7310
$fitnessValues[127,0] = 7310
$fitnessValues[128,0] = 7310
$fitnessValues[129,0] = 7310
$fitnessValues[130,0] = 7310
$fitnessValues[131,0] = 7310
$fitnessValues[132,0] = 7310
$fitnessValues[133,0] = 7310
$fitnessValues[134,0] = 7310
$fitnessValues[135,0] = 7310
$fitnessValues[136,0] = 7310
$fitnessValues[137,0] = 7310
$fitnessValues[138,0] = 7310
$fitnessValues[139,0] = 7310
$fitnessValues[140,0] = 7295
$fitnessValues[141,0] = 7295
$fitnessValues[142,0] = 7295
$fitnessValues[143,0] = 7295
$fitnessValues[144,0] = 7295

$ws.Range("C2:C146").Value = $fitnessValues

Write-Output "Updated fitness values for rows 2-146"
